$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 873.55554
$ws.Range("I18").Value = 674.9286
$ws.Range("K18").Value = 674.9286
$ws.Range("M18").Value = -390.9286

$ws.Range("H32").Value = 2124
$ws.Range("I32").Value = 2165.3333
$ws.Range("K32").Value = 2165.3333
$ws.Range("M32").Value = -1839.3333

$ws.Range("H62").Value = 25862.25
$ws.Range("J62").Value = 32416.334
$ws.Range("L62").Value = 32416.334
$ws.Range("N62").Value = -33664.334

$ws.Range("H65").Value = 25862.25
$ws.Range("J65").Value = 32416.334
$ws.Range("L65").Value = 162081.67
$ws.Range("N65").Value = -168321.67

$ws.Range("H116").Value = 5528.4707
$ws.Range("I116").Value = 5202.2666
$ws.Range("J116").Value = 7975
$ws.Range("K116").Value = 5202.2666
$ws.Range("L116").Value = 7975
$ws.Range("M116").Value = -1760.2666
$ws.Range("N116").Value = -14859

$ws.Range("H137").Value = 3190.8975
$ws.Range("I137").Value = 2463.2307
$ws.Range("K137").Value = 7389.6921
$ws.Range("M137").Value = -4839.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6026137.5
$ws.Range("I32").Value = 6099139
$ws.Range("K32").Value = 6099139
$ws.Range("M32").Value = -6098852

$ws.Range("H45").Value = 1912.8
$ws.Range("I45").Value = 1845.8334
$ws.Range("J45").Value = 2013.25
$ws.Range("K45").Value = 1845.8334
$ws.Range("L45").Value = 2013.25
$ws.Range("M45").Value = -1468.8334
$ws.Range("N45").Value = -2767.25

$ws.Range("H74").Value = 10426812
$ws.Range("I74").Value = 13892072
$ws.Range("K74").Value = 13892072
$ws.Range("M74").Value = -13891198

$ws.Range("H77").Value = 10426812
$ws.Range("I77").Value = 13892072
$ws.Range("K77").Value = 69460360
$ws.Range("M77").Value = -69455992

$ws.Range("H97").Value = 1361
$ws.Range("I97").Value = 1490.7646
$ws.Range("K97").Value = 1490.7646
$ws.Range("M97").Value = -994.7646

$ws.Range("H132").Value = 4839.4814
$ws.Range("I132").Value = 416.57144
$ws.Range("J132").Value = 9602.615
$ws.Range("K132").Value = 1249.71432
$ws.Range("L132").Value = 28807.845
$ws.Range("M132").Value = 1280.28568
$ws.Range("N132").Value = -33867.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2100.5208
$ws.Range("I20").Value = 2563.861
$ws.Range("K20").Value = 2563.861
$ws.Range("M20").Value = -2316.861

$ws.Range("H133").Value = 37396.8
$ws.Range("J133").Value = 47333
$ws.Range("L133").Value = 47333
$ws.Range("N133").Value = -57453

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2580.9092
$ws.Range("I62").Value = 2510.111
$ws.Range("J62").Value = 2899.5
$ws.Range("K62").Value = 2510.111
$ws.Range("L62").Value = 2899.5
$ws.Range("M62").Value = -1886.111
$ws.Range("N62").Value = -4147.5

$ws.Range("H65").Value = 2580.9092
$ws.Range("I65").Value = 2510.111
$ws.Range("J65").Value = 2899.5
$ws.Range("K65").Value = 12550.555
$ws.Range("L65").Value = 14497.5
$ws.Range("M65").Value = -9430.555
$ws.Range("N65").Value = -20737.5

$ws.Range("H105").Value = 1830.4667
$ws.Range("I105").Value = 1830.4667
$ws.Range("K105").Value = 1830.4667
$ws.Range("M105").Value = -83.46669999999995

$ws.Range("H132").Value = 1829.7097
$ws.Range("I132").Value = 1500.6923
$ws.Range("J132").Value = 3540.6
$ws.Range("K132").Value = 4502.0769
$ws.Range("L132").Value = 10621.8
$ws.Range("M132").Value = -1972.0769
$ws.Range("N132").Value = -15681.8

$ws.Range("H134").Value = 627180
$ws.Range("I134").Value = 1251973
$ws.Range("K134").Value = 3755919
$ws.Range("M134").Value = -3753384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 5597.2
$ws.Range("I33").Value = 4662
$ws.Range("K33").Value = 27972
$ws.Range("M33").Value = -27689

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1088.1177
$ws.Range("I97").Value = 1262
$ws.Range("K97").Value = 1262
$ws.Range("M97").Value = -766

$ws.Range("H101").Value = 120996
$ws.Range("J101").Value = 120996
$ws.Range("L101").Value = 120996
$ws.Range("N101").Value = -127486

$ws.Range("H102").Value = 2615.547
$ws.Range("I102").Value = 2301.7727
$ws.Range("K102").Value = 2301.7727
$ws.Range("M102").Value = -679.7727

$ws.Range("H132").Value = 26320468
$ws.Range("I132").Value = 40001428
$ws.Range("K132").Value = 120004284
$ws.Range("M132").Value = -120001754

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4695.357
$ws.Range("I40").Value = 3970.5557
$ws.Range("K40").Value = 3970.5557
$ws.Range("M40").Value = -3834.5557

$ws.Range("H46").Value = 6103.5835
$ws.Range("J46").Value = 14500
$ws.Range("L46").Value = 14500
$ws.Range("N46").Value = -14876

$ws.Range("H122").Value = 4681.9736
$ws.Range("I122").Value = 4830.933
$ws.Range("K122").Value = 14492.799
$ws.Range("M122").Value = -12042.799

$ws.Range("H132").Value = 168328.86
$ws.Range("I132").Value = 4982.2383
$ws.Range("K132").Value = 14946.7149
$ws.Range("M132").Value = -12416.7149

$ws.Range("H136").Value = 43879.9
$ws.Range("I136").Value = 6760.3887
$ws.Range("K136").Value = 20281.1661
$ws.Range("M136").Value = -17731.1661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6700
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H108").Value = 61825
$ws.Range("J108").Value = 42062.5
$ws.Range("L108").Value = 42062.5
$ws.Range("N108").Value = -49742.5

$ws.Range("H123").Value = 44985.715
$ws.Range("J123").Value = 44985.715
$ws.Range("L123").Value = 44985.715
$ws.Range("N123").Value = -54785.715

$ws.Range("H126").Value = 4133.758
$ws.Range("I126").Value = 4442.04
$ws.Range("J126").Value = 3170.375
$ws.Range("K126").Value = 13326.12
$ws.Range("L126").Value = 9511.125
$ws.Range("M126").Value = -10856.12
$ws.Range("N126").Value = -14451.125

$ws.Range("H132").Value = 2862.8948
$ws.Range("I132").Value = 2862.8948
$ws.Range("K132").Value = 8588.6844
$ws.Range("M132").Value = -6058.6844

$ws.Range("H136").Value = 1318.1875
$ws.Range("I136").Value = 616.0833
$ws.Range("K136").Value = 1848.2499
$ws.Range("M136").Value = 701.7501

$ws.Range("H137").Value = 80141.60000000001
$ws.Range("J137").Value = 80141.60000000001
$ws.Range("L137").Value = 80141.60000000001
$ws.Range("N137").Value = -90341.60000000001
